$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) for rows 2-9
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 280
$ws1.Range("F3").Value = 176
$ws1.Range("F4").Value = 2110
$ws1.Range("F5").Value = 1665
$ws1.Range("F6").Value = 304
$ws1.Range("F7").Value = 88
$ws1.Range("F8").Value = 711
$ws1.Range("F9").Value = 156

# Sheet "演出" - update column F (想去人数) for row 2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 21

# Sheet "全部类型" - update column F (想去人数) for rows 2-10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 280
$ws4.Range("F3").Value = 176
$ws4.Range("F4").Value = 2110
$ws4.Range("F5").Value = 1665
$ws4.Range("F6").Value = 304
$ws4.Range("F7").Value = 21
$ws4.Range("F8").Value = 88
$ws4.Range("F9").Value = 711
$ws4.Range("F10").Value = 156
